$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

# Set the B column (Key) values first for rows 20-22
$ws.Range("B20").Value = "VesselPolicyNumber"
$ws.Range("B21").Value = "VesselClaimNumber"
$ws.Range("B22").Value = "VesselName"

# Shared description value used across the new rows
$ws.Range("D20").Value = "From VesselRiskinfo excel"

# Set the C column (Value) values for rows 20-22
$ws.Range("C20").Value = "policy"
$ws.Range("C21").Value = "claim"
$ws.Range("C22").Value = "Vessel Name"

# Fill in remaining D column values reusing the same shared string
$ws.Range("D21").Value = "From VesselRiskinfo excel"
$ws.Range("D22").Value = "From VesselRiskinfo excel"

$ws.Range("B22").Select()
